$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S (2022 data) mirrors the existing formatting used by column Q
# (the "last data column" formatting previously applied to R is now shared by
# both R and the new S column).
$ws.Range("Q4:Q7").Copy() | Out-Null
$ws.Range("R4:R7").PasteSpecial(-4122) | Out-Null
$ws.Range("S4:S7").PasteSpecial(-4122) | Out-Null

# 2021 values in R stay the same - only formatting moved to match Q.
$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 49.9
$ws.Range("R6").Value = 33.6
$ws.Range("R7").Value = 25.1

# New 2022 figures in column S.
$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 49.7
$ws.Range("S6").Value = 34.9
$ws.Range("S7").Value = 21

# Move the active selection the way the author left it.
$ws.Range("R12").Select() | Out-Null
